# Applies the "test remove pawn completed" edit to EventResult_Data.xlsx
# Row 7  (Id=105, Event5): descrption/condition/option updated to reflect a removed pawn.
# Row 8  (Id=106, Event6): descrption/condition/option updated to reflect an added pawn.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 -> Id 105
$ws.Range("C7").Value = "remove 1003 character"
$ws.Range("D7").Value = 21
$ws.Range("E7").Value = "[106]"

# Row 8 -> Id 106
$ws.Range("C8").Value = "add 1003 character"
$ws.Range("D8").Value = 22
$ws.Range("E8").Value = "[107]"
